$d = $word.ActiveDocument

$replacements = @(
    @{old="2023-11-17 Friday"; new="2023-11-18 Saturday"},
    @{old="93×48=4464"; new="27×67=1809"},
    @{old="21×26=546"; new="85×77=6545"},
    @{old="93×24=2232"; new="41×21=861"},
    @{old="12×88=1056"; new="40×66=2640"},
    @{old="72×24=1728"; new="95×75=7125"},
    @{old="17×81=1377"; new="74×71=5254"},
    @{old="29×99=2871"; new="44×53=2332"},
    @{old="68×97=6596"; new="77×62=4774"},
    @{old="28×39=1092"; new="35×45=1575"},
    @{old="33×34=1122"; new="80×73=5840"},
    @{old="77×96=7392"; new="76×15=1140"},
    @{old="67×96=6432"; new="48×44=2112"},
    @{old="19×54=1026"; new="29×75=2175"},
    @{old="80×92=7360"; new="81×86=6966"},
    @{old="60×37=2220"; new="24×59=1416"},
    @{old="97×70=6790"; new="89×80=7120"},
    @{old="43×95=4085"; new="45×64=2880"},
    @{old="59×39=2301"; new="86×57=4902"},
    @{old="83×13=1079"; new="20×92=1840"},
    @{old="93×29=2697"; new="49×37=1813"},
    @{old="84×29=2436"; new="52×92=4784"},
    @{old="52×44=2288"; new="80×91=7280"},
    @{old="80×89=7120"; new="79×52=4108"},
    @{old="85×66=5610"; new="70×94=6580"},
    @{old="59×51=3009"; new="29×55=1595"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

Write-Output "done"
